$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in D7 and E7 (Button / Timer columns for Step 5)
$ws.Range("D7").Value = "T"
$ws.Range("E7").Value = "X"

# Update the active selection to F7
$ws.Range("F7").Select()
